$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 288, shifting rows 288:307 down to 289:308
$ws.Rows("288:288").Insert()

# Populate the newly inserted row 288 with the new record
$ws.Range("A288").Value = 10
$ws.Range("B288").Value = "Vega Modelo de Temuco"
$ws.Range("C288").Value = "La Araucanía"
$ws.Range("D288").Value = 45223
$ws.Range("E288").Value = 9
$ws.Range("F288").Value = "Fruta"
$ws.Range("G288").Value = 100104
$ws.Range("H288").Value = "Frutos de pepita"
$ws.Range("I288").Value = 100104001
$ws.Range("J288").Value = "Granada"
$ws.Range("K288").Value = "Wonderfull"
$ws.Range("L288").Value = "Primera"
$ws.Range("M288").Value = 80
$ws.Range("N288").Value = 17000
$ws.Range("O288").Value = 17000
$ws.Range("P288").Value = 17000
$ws.Range("Q288").Value = "$/bandeja 10 kilos granel"
$ws.Range("R288").Value = "Provincia de Limarí"
$ws.Range("S288").Value = 1700
$ws.Range("T288").Value = 10
